# Updates cryptos list prices/volumes/coins per the upstream scrape refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = '@'
    $rng.Value = $text
    $rng.ClearFormats()
}

Set-TextValue 'D2' '27.645.01'
Set-TextValue 'E2' '  -0.88%  '
Set-TextValue 'D3' '1.633.97'
Set-TextValue 'E3' '  +0.05%  '
Set-TextValue 'D4' '1.00'
Set-TextValue 'E4' '  -0.23%  '
Set-TextValue 'D5' '211.61'
Set-TextValue 'E5' '  -0.46%  '
Set-TextValue 'D7' '1.00'
Set-TextValue 'E7' '  -0.23%  '
Set-TextValue 'D8' '23.20'
Set-TextValue 'E8' '  +0.39%  '
Set-TextValue 'E9' '  +0.22%  '
Set-TextValue 'E10' '  -0.13%  '
Set-TextValue 'E11' '  -3.20%  '
Set-TextValue 'D12' '1.865.26'
Set-TextValue 'E12' '  +0.04%  '
Set-TextValue 'D13' '1.629.94'
Set-TextValue 'E13' '  -0.03%  '
Set-TextValue 'D14' '4.04'
Set-TextValue 'E14' '  -0.31%  '
Set-TextValue 'D15' '0.561'
Set-TextValue 'E15' '  +0.83%  '
Set-TextValue 'D16' '65.18'
Set-TextValue 'E16' '  +0.95%  '
Set-TextValue 'D17' '27.634.47'
Set-TextValue 'E17' '  -0.76%  '
Set-TextValue 'D18' '229.88'
Set-TextValue 'E18' '  -0.74%  '
Set-TextValue 'E19' '  -0.44%  '
Set-TextValue 'D20' '7.61'
Set-TextValue 'E20' '  -0.12%  '
Set-TextValue 'E21' '  -0.32%  '
Set-TextValue 'D22' '10.67'
Set-TextValue 'E22' '  +7.05%  '
Set-TextValue 'D23' '4.38'
Set-TextValue 'E23' '  +1.71%  '
Set-TextValue 'E24' '  +2.80%  '
Set-TextValue 'D25' '149.59'
Set-TextValue 'E25' '  -0.19%  '
Set-TextValue 'E26' '  -0.74%  '
Set-TextValue 'E27' '  -0.77%  '
Set-TextValue 'D28' '15.63'
Set-TextValue 'E28' '  +0.01%  '
Set-TextValue 'D29' '1.00'
Set-TextValue 'E29' '  -0.23%  '
Set-TextValue 'E30' '  -0.31%  '
Set-TextValue 'D31' '0.0481'
Set-TextValue 'E31' '  -0.55%  '
Set-TextValue 'D32' '3.28'
Set-TextValue 'E32' '  -0.76%  '
Set-TextValue 'D33' '1.466.89'
Set-TextValue 'E33' '  -0.25%  '
Set-TextValue 'E34' '  -0.03%  '
Set-TextValue 'E35' '  +0.09%  '
Set-TextValue 'E36' '  -1.63%  '
Set-TextValue 'E37' '  +0.16%  '
Set-TextValue 'B38' 'TrustWalletToken'
Set-TextValue 'C38' 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue 'D38' '0.924'
Set-TextValue 'E38' '  -0.36%  '
Set-TextValue 'E39' '  -0.39%  '
Set-TextValue 'B40' 'ImmutableX'
Set-TextValue 'C40' 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue 'D40' '0.557'
Set-TextValue 'E40' '  -1.70%  '
Set-TextValue 'D41' '68.88'
Set-TextValue 'E41' '  -0.80%  '
Set-TextValue 'E42' '  -0.25%  '
Set-TextValue 'E43' '  -0.45%  '
Set-TextValue 'E44' '  +0.17%  '
Set-TextValue 'E45' '  -0.80%  '
Set-TextValue 'D46' '5.36'
Set-TextValue 'E46' '  -0.85%  '
Set-TextValue 'D47' '1.774.67'
Set-TextValue 'E47' '  -0.06%  '
Set-TextValue 'E48' '  +2.72%  '
Set-TextValue 'D49' '87.66'
Set-TextValue 'E49' '  +2.06%  '
Set-TextValue 'D50' '0.0₆0105'
Set-TextValue 'E50' '  -0.70%  '
Set-TextValue 'D51' '0.0998'
Set-TextValue 'E51' '  +0.64%  '
